$wb = $excel.ActiveWorkbook

# zh-cn sheet: update the Correspond Handoff Datetime / Correspond Handback DateTime
# for the f13b4c6d-...zh-cn.xlf row (row 5 of the table).
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D5").Value = "2016-01-21 02:56:59"
$wsZhCn.Range("G5").Value = "2016-01-21 02:57:47"

# de-de sheet: update the Correspond Handoff Datetime / Correspond Handback DateTime
# for the f13b4c6d-...de-de.xlf row (row 5 of the table).
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D5").Value = "2016-01-21 02:57:12"
$wsDeDe.Range("G5").Value = "2016-01-21 02:58:08"
